# Generate Report for Handoff
# Update the handoff timestamps for file 4bf3d61f-5567-4248-b1d8-38d08e9bc807
# across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-27-14 08:27:20"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-14 08:27:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-14 08:27:20"
